# Pawsz.xlsx update:
#  - Reorder prey/age classes so "a" (adult) comes before "aa" (aged adult) etc,
#    make UNID ("u") sex the 1st-row-anchored block stays last,
#    add a new "aa" (aged adult) age class,
#    and extend the "u" (unknown sex) average rows to cover it too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Sex), Column B (Ageclass), Column C (Pawidth) ------------
# Row layout after the edit (rows 2-16):
#   f / j   / 40
#   f / sa  / 44
#   f / a   / 47
#   f / aa  / 47
#   f / u   / 46
#   m / j   / 42
#   m / sa  / 45
#   m / a   / 49
#   m / a   / 49
#   m / u   / 48
#   u / j   / =(C2+C7)/2
#   u / sa  / =(C3+C8)/2
#   u / a   / =(C4+C9)/2
#   u / aa  / 48
#   u / u   / =(C6+C11)/2

$rows = @(
    @{ A = "f"; B = "j";  C = 40 },
    @{ A = "f"; B = "sa"; C = 44 },
    @{ A = "f"; B = "a";  C = 47 },
    @{ A = "f"; B = "aa"; C = 47 },
    @{ A = "f"; B = "u";  C = 46 },
    @{ A = "m"; B = "j";  C = 42 },
    @{ A = "m"; B = "sa"; C = 45 },
    @{ A = "m"; B = "a";  C = 49 },
    @{ A = "m"; B = "a";  C = 49 },
    @{ A = "m"; B = "u";  C = 48 },
    @{ A = "u"; B = "j";  C = '=(C2+C7)/2' },
    @{ A = "u"; B = "sa"; C = '=(C3+C8)/2' },
    @{ A = "u"; B = "a";  C = '=(C4+C9)/2' },
    @{ A = "u"; B = "aa"; C = 48 },
    @{ A = "u"; B = "u";  C = '=(C6+C11)/2' }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    if ($row.C -is [string] -and $row.C.StartsWith("=")) {
        $ws.Cells.Item($r, 3).Formula = $row.C
    } else {
        $ws.Cells.Item($r, 3).Value = $row.C
    }
    $r++
}

# New used range is A1:C16 - select it (matches the saved sqref, no activeCell override)
$ws.Range("A1:C16").Select()
